$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162, pushing the existing rows 162-191 down to 163-192
# (this mirrors the OOXML diff, which shows every row from 163..192 now
# holding what used to be in the row above it, and a brand-new row 162
# with a week's worth of fresh data).
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly record.
$ws.Cells.Item(162, 1).Value = 4
$ws.Cells.Item(162, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(162, 3).Value = "Los Lagos"
$ws.Cells.Item(162, 4).Value = 44522
$ws.Cells.Item(162, 5).Value = 10
$ws.Cells.Item(162, 6).Value = 100112037
$ws.Cells.Item(162, 7).Value = "Cebollín"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 60
$ws.Cells.Item(162, 11).Value = 5000
$ws.Cells.Item(162, 12).Value = 5000
$ws.Cells.Item(162, 13).Value = 5000
$ws.Cells.Item(162, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(162, 15).Value = "Región Metropolitana"
$ws.Cells.Item(162, 16).Value = 139
$ws.Cells.Item(162, 17).Value = 36
$ws.Cells.Item(162, 18).Value = "Hortaliza"
